$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row data (rows 2..11): Colaborador_id, Colaborador_nome, Departamento,
# Motivo_da_ausência, Horas_de_ausência, Data_da_ausência (serial), Salário
$data = @(
    @(82899, "Lorena Moraes", "Operações", "Viagem de negócios", 8, 45087, 10771.69),
    @(50613, "Maria Eduarda Rocha", "Operações", "Outros", 6, 45093, 8245.559999999999),
    @(55990, "Vitor Hugo Pereira", "Operações", "Outros", 4, 45100, 9140.18),
    @(26547, "Daniela Ribeiro", "Operações", "Consulta médica", 8, 45098, 8917.700000000001),
    @(30677, "Luna da Rocha", "Recursos Humanos", "Outros", 6, 45090, 6022.75),
    @(83966, "Lucca Costa", "Jurídico", "Problemas pessoais", 7, 45094, 11875.26),
    @(94188, "Rafael da Rosa", "Operações", "Viagem de negócios", 3, 45092, 9216.9),
    @(50584, "Breno Pereira", "Vendas", "Outros", 5, 45080, 6947.44),
    @(65257, "Isadora Martins", "Engenharia", "Problemas pessoais", 4, 45096, 11036.01),
    @(21992, "Sofia Ramos", "Engenharia", "Consulta médica", 1, 45082, 9672.379999999999)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r++
}
